# Gameflow workbook edit:
# Add two new ALLY rows ("Illyia" and "Dahlia") to the BATTLE scene block,
# right after the MUSIC row (old row 11), pushing the existing ENEMY rows
# (and everything below them) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 12 (shifts old row 12 "ENEMY/Slime"
# and everything after it down to rows 14+).
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Fill in the new ALLY rows (write order matches the shared-string table
# order produced by the original edit: Illyia, then ALLY, then Dahlia).
$ws.Range("B12").Value = "Illyia"
$ws.Range("A12").Value = "ALLY"

$ws.Range("B13").Value = "Dahlia"
$ws.Range("A13").Value = "ALLY"

# Match the author's final selection (active cell on the newly added B13).
$ws.Range("B13").Select()
